$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "NSE:ASTRAL"
$ws.Range("C2").Value = "NSE:ADANIPORTS"
$ws.Range("D2").Value = ""
$ws.Range("F2").Value = "NSE:ASTRAL"
$ws.Range("B3").Value = "NSE:GOKULAGRO"
$ws.Range("C3").Value = "NSE:ADROITINFO"
$ws.Range("E3").Value = "NSE:MOTHERSON"
$ws.Range("F3").Value = "NSE:MCX"
$ws.Range("B4").Value = "NSE:GOODLUCK"
$ws.Range("C4").Value = "NSE:ALANKIT"
$ws.Range("E4").Value = "NSE:PPLPHARMA"
$ws.Range("F4").Value = "NSE:MUTHOOTFIN"
$ws.Range("B5").Value = "NSE:MCX"
$ws.Range("C5").Value = "NSE:ALOKINDS"
$ws.Range("B6").Value = "NSE:MOKSH"
$ws.Range("C6").Value = "NSE:AMBIKCO"
$ws.Range("B7").Value = "NSE:RAMRAT"
$ws.Range("C7").Value = "NSE:ASMS"
$ws.Range("C8").Value = "NSE:BALKRISIND"
$ws.Range("C9").Value = "NSE:BASML"
$ws.Range("C10").Value = "NSE:BCONCEPTS"
$ws.Range("C11").Value = "NSE:BLS"
$ws.Range("C12").Value = "NSE:CANBK"
$ws.Range("C13").Value = "NSE:DATAMATICS"
$ws.Range("C14").Value = "NSE:DREAMFOLKS"
$ws.Range("C15").Value = "NSE:DWARKESH"
$ws.Range("C16").Value = "NSE:ELDEHSG"
$ws.Range("C17").Value = "NSE:ELECTCAST"
$ws.Range("C18").Value = "NSE:FINEORG"
$ws.Range("C19").Value = "NSE:GARFIBRES"
$ws.Range("C20").Value = "NSE:GEECEE"
$ws.Range("C21").Value = "NSE:GENUSPAPER"
$ws.Range("C22").Value = "NSE:GTLINFRA"
$ws.Range("C23").Value = "NSE:HFCL"
$ws.Range("C24").Value = "NSE:HITECH"
$ws.Range("C25").Value = "NSE:IGPL"
$ws.Range("C26").Value = "NSE:INDIGOPNTS"
$ws.Range("C27").Value = "NSE:INDOAMIN"
$ws.Range("C28").Value = "NSE:INDOWIND"
$ws.Range("C29").Value = "NSE:INGERRAND"
$ws.Range("C30").Value = "NSE:IRB"
$ws.Range("C31").Value = "NSE:JISLJALEQS"
$ws.Range("C32").Value = "NSE:KALAMANDIR"
$ws.Range("C33").Value = "NSE:KIOCL"
$ws.Range("C34").Value = "NSE:LINCOLN"
$ws.Range("C35").Value = "NSE:MAITHANALL"
$ws.Range("C36").Value = "NSE:MOTILALOFS"
$ws.Range("C37").Value = "NSE:MTNL"
$ws.Range("C38").Value = "NSE:NGIL"
$ws.Range("C39").Value = "NSE:NIRAJ"
$ws.Range("C40").Value = "NSE:NLCINDIA"
$ws.Range("C41").Value = "NSE:OLECTRA"
$ws.Range("C42").Value = "NSE:ONWARDTEC"
$ws.Range("C43").Value = "NSE:PENIND"
$ws.Range("C44").Value = "NSE:PREMIERPOL"
$ws.Range("C45").Value = "NSE:PRITI"
$ws.Range("C46").Value = "NSE:RECLTD"
$ws.Range("C47").Value = "NSE:REFEX"
$ws.Range("C48").Value = "NSE:RGL"
$ws.Range("C49").Value = "NSE:RICOAUTO"

# Delete now-unused trailing rows 50-69 (data previously ran to row 69,
# the updated table only spans rows 1-49)
$ws.Rows("50:69").Delete()
